$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = "2012-04-23"
    $ws.Cells.Item($r, 9).Value = "許忠信"
    $ws.Cells.Item($r, 10).Value = 1749
}
